# Updated capital structure database
# Applies the updated financial figures (columns D..AK) for the three
# Belgium / Bank (Money Center) rows (2, 3, 4) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.005105
$ws.Range("E2").Value = -0.00425
$ws.Range("F2").Value = 0.0268
$ws.Range("K2").Value = 4096
$ws.Range("L2").Value = 0.2592487104022279
$ws.Range("M2").Value = 3681.4
$ws.Range("N2").Value = 0.08372640973040342
$ws.Range("O2").Value = 0.8987792968749999
$ws.Range("P2").Value = 590.7
$ws.Range("Q2").Value = 0.01343434297488708
$ws.Range("R2").Value = 0.1442138671875
$ws.Range("S2").Value = 3090.7
$ws.Range("T2").Value = 0.8395447384147335
$ws.Range("U2").Value = 73010.8
$ws.Range("V2").Value = 1.660491159761107
$ws.Range("W2").Value = 0.08848224394917192
$ws.Range("X2").Value = 0.1078940081287463
$ws.Range("Y2").Value = -0.01941176417957434
$ws.Range("Z2").Value = 0.1333669293571848
$ws.Range("AB2").Value = 0.03620400601274515
$ws.Range("AC2").Value = -0.03620400601274515
$ws.Range("AD2").Value = 138835.7
$ws.Range("AF2").Value = 138835.7
$ws.Range("AG2").Value = 65824.90000000001
$ws.Range("AH2").Value = 0.7594738877635252
$ws.Range("AI2").Value = 0.7092689052276229
$ws.Range("AJ2").Value = 0.5995293016121966
$ws.Range("AK2").Value = 0.5363216386657325

# Row 3
$ws.Range("D3").Value = 0.009389999999999999
$ws.Range("E3").Value = 0.057
$ws.Range("K3").Value = 2212.5
$ws.Range("L3").Value = 0.26482099895867
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 39865
$ws.Range("V3").Value = 2.694436746804729
$ws.Range("W3").Value = 0.08876665503171527
$ws.Range("X3").Value = 0.1573012779342385
$ws.Range("Y3").Value = -0.06853462290252324
$ws.Range("Z3").Value = 0.1135704964384753
$ws.Range("AB3").Value = 0.03586778251236599
$ws.Range("AC3").Value = -0.03586778251236599
$ws.Range("AD3").Value = 100445.6
$ws.Range("AF3").Value = 100445.6
$ws.Range("AG3").Value = 60580.60000000001
$ws.Range("AH3").Value = 0.8716141578207043
$ws.Range("AI3").Value = 0.7561886722567105
$ws.Range("AJ3").Value = 0.8037131231600552
$ws.Range("AK3").Value = 0.6516397322043233

# Row 4
$ws.Range("D4").Value = -0.0196
$ws.Range("E4").Value = -0.0655
$ws.Range("F4").Value = 0.0268
$ws.Range("K4").Value = 1883.5
$ws.Range("L4").Value = 0.2529953793251666
$ws.Range("M4").Value = 3681.4
$ws.Range("N4").Value = 0.1261872688446259
$ws.Range("O4").Value = 1.954552694451818
$ws.Range("P4").Value = 590.7
$ws.Range("Q4").Value = 0.0202474112311948
$ws.Range("R4").Value = 0.3136182638704539
$ws.Range("S4").Value = 3090.7
$ws.Range("T4").Value = 0.8395447384147335
$ws.Range("U4").Value = 33145.8
$ws.Range("V4").Value = 1.136137875718531
$ws.Range("W4").Value = 0.08819783286662858
$ws.Range("X4").Value = 0.05848673832325402
$ws.Range("Y4").Value = 0.02971109454337456
$ws.Range("Z4").Value = 0.1657996009122007
$ws.Range("AB4").Value = 0.03654022951312431
$ws.Range("AC4").Value = -0.03654022951312431
$ws.Range("AD4").Value = 38390.1
$ws.Range("AF4").Value = 38390.1
$ws.Range("AG4").Value = 5244.299999999996
$ws.Range("AH4").Value = 0.568201799177671
$ws.Range("AI4").Value = 0.6102054570250535
$ws.Range("AJ4").Value = 0.1523690816539989
$ws.Range("AK4").Value = 0.1761747671965491

# T3 previously held 0; the cell is now fully cleared (no <c> element)
$ws.Range("T3").ClearContents()

